# Import_ATR_ExistedOppt2.xlsx - "Add last source code"
#
# This script reproduces the meaningful content edits made to Sheet1:
#   - F2 (Exp Doc Nbr) value changes from the old document number to a new one
#   - Two new data points are filled in on row 2: the Opportunity ID (CY2)
#     and the Quote Start/End Date (DD2/DE2)
#   - The active selection is moved to the newly entered date range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F2: Exp Doc Nbr updated to the new document number ---
$ws.Range("F2").Value = "8819091909"

# --- CY2: Oppty ID (plain/default formatting, like the rest of the sparse row) ---
$ws.Range("CY2").Value = "16D7DEDB-4F7E-4CED-B4C1-97214D6BF9AE"
$ws.Range("CY2").Style = "Normal"

# --- DD2 / DE2: Quote Start Date / Quote End Date ---
$ws.Range("DD2").Value = 42382
$ws.Range("DE2").Value = 42656

# Give the new date cells the same date formatting already used elsewhere
# in this row/sheet (reuse existing style instead of inventing a new one).
$ws.Range("DN2").Copy() | Out-Null
$ws.Range("DD2:DE2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Reflect the place the user ended up after entering the new data.
$ws.Range("DD2:DE2").Select() | Out-Null
